$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new D value (price) and new E value (volume %)
# $null means the D cell is unchanged for that row.
$updates = @(
    @{ Row = 2;  D = "71.155.79";  E = "  -0.40%  " },
    @{ Row = 3;  D = "3.834.14";   E = "  +0.44%  " },
    @{ Row = 4;  D = $null;        E = "  -0.08%  " },
    @{ Row = 5;  D = "704.08";     E = "  +0.52%  " },
    @{ Row = 6;  D = "172.05";     E = "  -1.29%  " },
    @{ Row = 7;  D = "3.833.99";   E = "  +0.44%  " },
    @{ Row = 8;  D = $null;        E = "  -0.01%  " },
    @{ Row = 9;  D = "0.526";      E = "  -0.63%  " },
    @{ Row = 10; D = $null;        E = "  -0.33%  " },
    @{ Row = 11; D = "7.37";       E = "  -0.74%  " },
    @{ Row = 12; D = $null;        E = "  -0.46%  " },
    @{ Row = 13; D = "0.0000255";  E = "  -1.73%  " },
    @{ Row = 14; D = "36.74";      E = "  +0.48%  " },
    @{ Row = 15; D = "4.480.69";   E = "  +0.44%  " },
    @{ Row = 16; D = "3.798.56";   E = "  +0.12%  " },
    @{ Row = 17; D = "71.153.32";  E = "  -0.41%  " },
    @{ Row = 18; D = "7.24";       E = "  -0.07%  " },
    @{ Row = 19; D = $null;        E = "  +0.26%  " },
    @{ Row = 20; D = "17.40";      E = "  -2.08%  " },
    @{ Row = 21; D = "10.72";      E = "  -3.85%  " },
    @{ Row = 22; D = "495.15";     E = "  +1.78%  " },
    @{ Row = 23; D = "0.737";      E = "  +2.83%  " },
    @{ Row = 24; D = "85.34";      E = "  +0.62%  " },
    @{ Row = 25; D = "0.0000145";  E = "  +0.92%  " },
    @{ Row = 26; D = "10.63";      E = "  +0.72%  " },
    @{ Row = 27; D = $null;        E = "  -2.12%  " },
    @{ Row = 28; D = "2.09";       E = "  -2.92%  " },
    @{ Row = 29; D = "1.00";       E = "  -0.07%  " },
    @{ Row = 30; D = "3.08";       E = "  -2.49%  " },
    @{ Row = 31; D = "7.44";       E = "  -2.26%  " },
    @{ Row = 32; D = "2.24";       E = "  -3.62%  " },
    @{ Row = 33; D = "29.47";      E = "  -0.76%  " },
    @{ Row = 34; D = $null;        E = "  -3.42%  " },
    @{ Row = 35; D = "9.22";       E = "  -0.94%  " },
    @{ Row = 36; D = "3.795.30";   E = "  +0.73%  " },
    @{ Row = 37; D = "0.999";      E = "  -0.09%  " },
    @{ Row = 38; D = "0.103";      E = "  -0.85%  " },
    @{ Row = 39; D = "2.34";       E = "  -2.13%  " },
    @{ Row = 40; D = $null;        E = "  +3.48%  " },
    @{ Row = 41; D = "6.00";       E = "  -0.63%  " },
    @{ Row = 42; D = "3.32";       E = "  -2.94%  " },
    @{ Row = 43; D = $null;        E = "  +0.01%  " },
    @{ Row = 44; D = $null;        E = "  -0.05%  " },
    @{ Row = 45; D = "0.000313";   E = "  +1.15%  " },
    @{ Row = 46; D = "163.94";     E = "  +0.13%  " },
    @{ Row = 47; D = "428.21";     E = "  +2.92%  " },
    @{ Row = 48; D = "48.93";      E = "  +0.45%  " },
    @{ Row = 49; D = "8.76";       E = "  +0.64%  " },
    @{ Row = 50; D = $null;        E = "  +0.02%  " },
    @{ Row = 51; D = "0.297";      E = "  -1.65%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cellD = $ws.Range("D$($u.Row)")
        # These price strings look numeric (e.g. "704.08"); assigning them
        # directly through .Value would make Excel auto-convert to a
        # number, same as typing into a General-formatted cell. Force the
        # cell to Text first only when the string would actually parse as
        # a plain number, so the value round-trips as text like the
        # original inline string cells did.
        if ($u.D -match '^[+-]?[0-9]*\.?[0-9]+$') {
            $cellD.NumberFormat = "@"
        }
        $cellD.Value = $u.D
    }
    $ws.Range("E$($u.Row)").Value = $u.E
}
